# se corrijen algunas palabras de Culture
# Fix a few words/typos on the "Preguntas" sheet (CVF culture questionnaire data).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Preguntas")
$ws2 = $wb.Worksheets.Item("Prefijos")

# Replace "instituciones educativas" / "educational institutions" wording with "compañias" / "companies"
$ws1.Range("E24").Value = "…ganar participación en el mercado y competir con otras compañias."
$ws1.Range("F24").Value = "…gaining market share and competing with other companies."

# Fix typo: "Emfasis" -> "Enfasis"
$ws1.Range("E15").Value = "...compromiso con la innovación y el desarrollo. Enfasis en posicionarse en prácticas (modelos, tecnología) de punta. "

# Restore view/selection state
$ws2.Activate()
$ws2.Cells.Select()

$ws1.Activate()
$ws1.Range("E13").Select()
